# ---------------------------------------------------------------------------
# Applies the "updated on Oct 9" commit:
#   * inserts a new "Battery" sheet between "Cost" and "result.csv"
#   * appends a "detministic setup" block to the "Cost" sheet
#   * appends the same "detministic setup" block to the "result.csv" sheet
#   * tweaks a couple of sheet-view / selection / active-tab details
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write one row of values starting at (row, startCol) on $ws.
function Set-RowValues {
    param($ws, [int]$row, [int]$startCol, [object[]]$values)
    $col = $startCol
    foreach ($v in $values) {
        if ($null -ne $v) {
            $ws.Cells.Item($row, $col).Value = $v
        }
        $col++
    }
}

# ---------------------------------------------------------------------------
# 1. Create the "Battery" sheet as a copy of "Cost" (so it inherits the same
#    sheetFormatPr / namespaces / default row height) positioned right after
#    "Cost", then wipe its contents so we can build it up from scratch.
# ---------------------------------------------------------------------------
$costSheet = $wb.Worksheets.Item("Cost")
$costSheet.Copy([System.Reflection.Missing]::Value, $costSheet)

$battery = $wb.Worksheets.Item("Cost (2)")
$battery.Name = "Battery"
$battery.Cells.Clear()

# ---------------------------------------------------------------------------
# 2. Populate the Battery sheet
# ---------------------------------------------------------------------------

# Row 1: header (yellow fill + centered)
Set-RowValues $battery 1 1 @("BCap", "pen_25", "pen_50", "pen_75", "pen_100")

# Rows 2-4: data (yellow fill + centered)
Set-RowValues $battery 2 1 @(3,  24025.820426014223, 17860.145015957976, 11190.917434693136, 8615.2970136409349)
Set-RowValues $battery 3 1 @(15, 23505.235062796452, 17339.559809490325, 10657.0707664399,   8002.8492457905486)
Set-RowValues $battery 4 1 @(30, 22851.462678712902, 16685.787744748824, 10016.551889938706, 7283.7349684510646)

# Row 7: section title "detministic setup" (merged A7:B7)
$battery.Range("A7").Value = "detministic setup"
$battery.Range("A7:B7").Merge()

# Row 8: header for the deterministic block (no fill, centered)
Set-RowValues $battery 8 1 @("BCap", "pen_25", "pen_50", "pen_75", "pen_100")

# Rows 9-11: deterministic data (no fill, centered)
Set-RowValues $battery 9  1 @(3,  23700.11, 16002.47, 10752.28,           8121.33)
Set-RowValues $battery 10 1 @(15, 22459.47, 14685.72, 9092.3700000000008, 6094.08)
Set-RowValues $battery 11 1 @(30, 20912.47, 13101.62, 7160.02,            3750.06)

# Row 16: header for the % difference block (yellow fill + centered)
Set-RowValues $battery 16 1 @("BCap", "pen_25", "pen_50", "pen_75", "pen_100")

# Rows 17-19: % difference between the stochastic (rows 2-4) and the
# deterministic (rows 9-11) results.
$battery.Range("A17").Value = 3
$battery.Range("A18").Value = 15
$battery.Range("A19").Value = 30
foreach ($col in @("B", "C", "D", "E")) {
    $battery.Range("$col" + "17").Formula = "=($col" + "2-$col" + "9)/$col" + "9"
    $battery.Range("$col" + "18").Formula = "=($col" + "3-$col" + "10)/$col" + "10"
    $battery.Range("$col" + "19").Formula = "=($col" + "4-$col" + "11)/$col" + "11"
}

# ---------------------------------------------------------------------------
# Styling for the Battery sheet
# ---------------------------------------------------------------------------

# Yellow fill + centered: header row, the 3 stochastic data rows, row 6 (B:E)
# and the % header row.
$yellowCentered = $battery.Range("A1:E4,B6:E6,A16:E19")
$yellowCentered.Interior.Color = 65535
$yellowCentered.HorizontalAlignment = -4108
$yellowCentered.VerticalAlignment = -4108

# No fill, but centered: blank spacer rows + the deterministic block
$plainCentered = $battery.Range("A5:E5,A6:A6,C7:E7,A8:E11")
$plainCentered.HorizontalAlignment = -4108
$plainCentered.VerticalAlignment = -4108

# Row 7 title cells (A7:B7) - centered too (style used for the merged title)
$battery.Range("A7:B7").HorizontalAlignment = -4108
$battery.Range("A7:B7").VerticalAlignment = -4108

# Percentage number format on the delta cells
$battery.Range("B17:E19").NumberFormat = "0.00%"

# Page setup: portrait orientation (matches the authored sheet)
$battery.PageSetup.Orientation = 1

# View: this tab is the active / selected one, with A1:E19 selected
$battery.Range("A1:E19").Select()
$battery.Activate()

# ---------------------------------------------------------------------------
# 3. "Cost" sheet: append the deterministic-setup block (rows 32-36)
# ---------------------------------------------------------------------------
$costSheet.Range("A32").Value = "detministic setup"
Set-RowValues $costSheet 33 1 @("BCap", "pen_25", "pen_50", "pen_75", "pen_100")
Set-RowValues $costSheet 34 1 @(3,  23700.11, 16002.47, 10752.28,           8121.33)
Set-RowValues $costSheet 35 1 @(15, 22459.47, 14685.72, 9092.3700000000008, 6094.08)
Set-RowValues $costSheet 36 1 @(30, 20912.47, 13101.62, 7160.02,            3750.06)

$costSheet.Range("A26:E36").Select()

# ---------------------------------------------------------------------------
# 4. "result.csv" sheet: append the same deterministic-setup block
#    (rows 31-35), noting that rows 33 & 34 are missing their E value in the
#    authored workbook.
# ---------------------------------------------------------------------------
$resultCsv = $wb.Worksheets.Item("result.csv")
$resultCsv.Range("A31").Value = "detministic setup"
Set-RowValues $resultCsv 32 1 @("BCap", "pen_25", "pen_50", "pen_75", "pen_100")
Set-RowValues $resultCsv 33 1 @(3,  23700.11, 16002.47, 10752.28)
Set-RowValues $resultCsv 34 1 @(15, 22459.47, 14685.72, 9092.3700000000008)
Set-RowValues $resultCsv 35 1 @(30, 20912.47, 13101.62, 7160.02, 3750.06)

$resultCsv.Columns.Item(1).ColumnWidth = 14.45

$resultCsv.Range("A31:E35").Select()
$resultCsv.Range("A31").Select()

# ---------------------------------------------------------------------------
# 5. Sheet-view / active-tab housekeeping
# ---------------------------------------------------------------------------
# "Solar" loses the tabSelected flag (Battery gets it instead, set above).
$solar = $wb.Worksheets.Item("Solar")

# Make the Battery tab the active / selected sheet (also drives workbook.xml
# bookViews/activeTab).
$battery.Activate()
$battery.Range("A1:E19").Select()
